{"js": "// The document is a flat list of single-run paragraphs (style \"capstyle\")\n// forming a \"Label :  Value\" report. The edit:\n//   - rewrites several paragraphs' text in place,\n//   - removes one paragraph entirely\n//     (\"Lymphovascular Invasion :  Large vessel (venous), intramural\"),\n//   - inserts one new paragraph after the \"Type of Polyp...\" line\n//     (\"Treatment Effect :  No known presurgical therapy\"),\n//   - appends five new paragraphs at the end of the body.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 1) Plain text -> text replacements, matched by the paragraph's current\n//    (pre-edit) text so the script does not depend on fragile indices.\nconst replacements = [\n  [\"Procedure :  Transverse colectomy\", \"Procedure :  Right hemicolectomy\"],\n  [\"Tumor Site :  Ileocecal valve\", \"Tumor Site :  Cecum\"],\n  [\"Histologic Grade :  G3, poorly differentiated\", \"Histologic Grade :  G1, well differentiated\"],\n  [\"Tumor Size :  Greatest dimension in Centimeters (cm)\", \"Tumor Extent :  Invades lamina propria / muscularis mucosae (intramucosal carcinoma)\"],\n  [\"Multiple Primary Sites  :  Present\", \"Macroscopic Tumor Perforation :  Not identified\"],\n  [\"Tumor Extent :  Invades submucosa\", \"Lymphovascular Invasion :  Not identified\"],\n  [\"Perineural Invasion :  Present\", \"Perineural Invasion :  Not identified\"],\n  [\"+Tumor Bud Score :  Low (0\", \"+Tumor Bud Score :  Low (0-4)\"],\n  [\"+Type of Polyp in which Invasive Carcinoma Arose :  Tubular adenoma\", \"+Type of Polyp in which Invasive Carcinoma Arose :  None identified\"],\n  [\"Regional Lymph Node Status :  All regional lymph nodes negative for tumor\", \"Regional Lymph Node Status :  Regional lymph nodes present\"],\n  [\"Number of Lymph Nodes with Tumor :  Other \", \"Number of Lymph Nodes with Tumor :  Exact number \"],\n  [\"Tumor Deposits :  Present\", \"Tumor Deposits :  Not identified\"],\n  [\"pT Category :  pT1\", \"Number of Tumor Deposits :  Specify number\"],\n];\n\nlet polypParagraph = null;       // anchor for the new \"Treatment Effect\" paragraph\nlet lastParagraph = null;        // anchor for the paragraphs appended at the end\nlet vascularInvasionParagraph = null; // paragraph to delete entirely\n\nfor (const item of paragraphs.items) {\n  const text = item.text;\n  const hit = replacements.find(([from]) => from === text);\n  if (hit) {\n    item.insertText(hit[1], \"Replace\");\n  }\n  if (text === \"+Type of Polyp in which Invasive Carcinoma Arose :  Tubular adenoma\") {\n    polypParagraph = item;\n  }\n  if (text === \"Lymphovascular Invasion :  Large vessel (venous), intramural\") {\n    vascularInvasionParagraph = item;\n  }\n  if (text === \"pT Category :  pT1\") {\n    lastParagraph = item;\n  }\n}\nawait context.sync();\n\n// 2) Remove the paragraph that has no replacement in the new report.\nif (vascularInvasionParagraph) {\n  vascularInvasionParagraph.delete();\n}\n\n// 3) Insert the new \"Treatment Effect\" paragraph right after the (now\n//    retextted) \"Type of Polyp...\" paragraph.\nif (polypParagraph) {\n  const inserted = polypParagraph.insertParagraph(\n    \"Treatment Effect :  No known presurgical therapy\",\n    \"After\"\n  );\n  inserted.style = \"capstyle\";\n}\n\n// 4) Append the five brand-new paragraphs at the end of the report, after\n//    what used to be \"pT Category :  pT1\" (now retexted in place above).\nconst newTailParagraphs = [\n  \"Distant Site :  Non-regional lymph node(s)\",\n  \"pT Category :  pT0\",\n  \"pT4 :  pT4a\",\n  \"pN Category :  pN not assigned (no nodes submitted or found)\",\n  \"pN2 :  pN2a\",\n];\n\nif (lastParagraph) {\n  let anchor = lastParagraph;\n  for (const text of newTailParagraphs) {\n    const inserted = anchor.insertParagraph(text, \"After\");\n    inserted.style = \"capstyle\";\n    anchor = inserted;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document is a flat list of single-run paragraphs (style \"capstyle\")\n# forming a \"Label :  Value\" report. The edit:\n#   - rewrites several paragraphs' text in place,\n#   - removes one paragraph entirely\n#     (\"Lymphovascular Invasion :  Large vessel (venous), intramural\"),\n#   - inserts one new paragraph after the \"Type of Polyp...\" line\n#     (\"Treatment Effect :  No known presurgical therapy\"),\n#   - appends five new paragraphs at the end of the body.\n\n$d = $word.ActiveDocument\n\n# 1) Plain text -> text replacements, located with Find so the script does\n#    not depend on fragile paragraph indices.\n$replacements = @(\n    @(\"Procedure :  Transverse colectomy\", \"Procedure :  Right hemicolectomy\"),\n    @(\"Tumor Site :  Ileocecal valve\", \"Tumor Site :  Cecum\"),\n    @(\"Histologic Grade :  G3, poorly differentiated\", \"Histologic Grade :  G1, well differentiated\"),\n    @(\"Tumor Size :  Greatest dimension in Centimeters (cm)\", \"Tumor Extent :  Invades lamina propria / muscularis mucosae (intramucosal carcinoma)\"),\n    @(\"Multiple Primary Sites  :  Present\", \"Macroscopic Tumor Perforation :  Not identified\"),\n    @(\"Tumor Extent :  Invades submucosa\", \"Lymphovascular Invasion :  Not identified\"),\n    @(\"Perineural Invasion :  Present\", \"Perineural Invasion :  Not identified\"),\n    @(\"+Tumor Bud Score :  Low (0\", \"+Tumor Bud Score :  Low (0-4)\"),\n    @(\"+Type of Polyp in which Invasive Carcinoma Arose :  Tubular adenoma\", \"+Type of Polyp in which Invasive Carcinoma Arose :  None identified\"),\n    @(\"Regional Lymph Node Status :  All regional lymph nodes negative for tumor\", \"Regional Lymph Node Status :  Regional lymph nodes present\"),\n    @(\"Number of Lymph Nodes with Tumor :  Other \", \"Number of Lymph Nodes with Tumor :  Exact number \"),\n    @(\"Tumor Deposits :  Present\", \"Tumor Deposits :  Not identified\"),\n    @(\"pT Category :  pT1\", \"Number of Tumor Deposits :  Specify number\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null\n}\n\n# 2) Remove the paragraph that has no replacement in the new report.\n#    Walk the paragraphs back-to-front so deleting one never invalidates\n#    the index of the paragraph we're about to look at next.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.TrimEnd([char]13)\n    if ($text -eq \"Lymphovascular Invasion :  Large vessel (venous), intramural\") {\n        $p.Range.Delete()\n    }\n}\n\n# 3) Insert the new \"Treatment Effect\" paragraph right after the (now\n#    retexted) \"Type of Polyp...\" paragraph.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.TrimEnd([char]13)\n    if ($text -eq \"+Type of Polyp in which Invasive Carcinoma Arose :  None identified\") {\n        $p.Range.InsertParagraphAfter() | Out-Null\n        $newPara = $d.Paragraphs.Item($i + 1)\n        $newPara.Range.Text = \"Treatment Effect :  No known presurgical therapy\"\n        $newPara.Style = \"capstyle\"\n        break\n    }\n}\n\n# 4) Append the five brand-new paragraphs at the end of the report, after\n#    what used to be \"pT Category :  pT1\" (now retexted in place above).\n$newTailParagraphs = @(\n    \"Distant Site :  Non-regional lymph node(s)\",\n    \"pT Category :  pT0\",\n    \"pT4 :  pT4a\",\n    \"pN Category :  pN not assigned (no nodes submitted or found)\",\n    \"pN2 :  pN2a\"\n)\n\n$anchorIndex = $d.Paragraphs.Count\n$anchor = $d.Paragraphs.Item($anchorIndex)\nforeach ($text in $newTailParagraphs) {\n    $anchor.Range.InsertParagraphAfter() | Out-Null\n    $anchorIndex = $anchorIndex + 1\n    $anchor = $d.Paragraphs.Item($anchorIndex)\n    $anchor.Range.Text = $text\n    $anchor.Style = \"capstyle\"\n}\n"}
